$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was inserted at row 6 (pushing the previous
# rows 6-24 down to rows 7-25).
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Macroferia Regional de Talca"
$ws.Range("C6").Value = "Maule"
$ws.Range("D6").Value = 45099
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100104
$ws.Range("H6").Value = "Frutos de pepita"
$ws.Range("I6").Value = 100104001
$ws.Range("J6").Value = "Granada"
$ws.Range("K6").Value = "Wonderfull"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 17000
$ws.Range("P6").Value = 17000
$ws.Range("Q6").Value = "$/caja 18 kilos granel"
$ws.Range("R6").Value = "Provincia de Limarí"
$ws.Range("S6").Value = 944
$ws.Range("T6").Value = 18
